$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.0630466474718978
$ws.Range("J2").Value = 0.06304664747189781
$ws.Range("M2").Value = 9.071155000000001
$ws.Range("N2").Value = 27.213465
$ws.Range("O2").Value = 0.1321733179750106
$ws.Range("P2").Value = 0.135675019324872
$ws.Range("Q2").Value = 2.588248466403333
$ws.Range("R2").Value = 23.29423619763
$ws.Range("S2").Value = 0.008333084583561546
$ws.Range("T2").Value = 0.008553855114118126
$ws.Range("I3").Value = 0.0630466474718978
$ws.Range("J3").Value = 0.06304664747189781
$ws.Range("N3").Value = 61.82726199999999
$ws.Range("O3").Value = 0.3002893736556623
$ws.Range("P3").Value = 0.3082450164524775
$ws.Range("Q3").Value = 5.880335931253777
$ws.Range("R3").Value = 52.92302338128399
$ws.Range("S3").Value = 0.01893223828042554
$ws.Range("T3").Value = 0.01943381488724869
$ws.Range("I4").Value = 0.0630466474718978
$ws.Range("J4").Value = 0.06304664747189781
$ws.Range("M4").Value = 17.81090666666666
$ws.Range("N4").Value = 53.43272
$ws.Range("O4").Value = 0.2595178486396241
$ws.Range("P4").Value = 0.2663933210482557
$ws.Range("Q4").Value = 5.081938503448888
$ws.Range("R4").Value = 45.73744653104
$ws.Range("S4").Value = 0.01636173031584771
$ws.Range("T4").Value = 0.01679520580099747
$ws.Range("I5").Value = 0.0630466474718978
$ws.Range("J5").Value = 0.06304664747189781
$ws.Range("M5").Value = 5.313972
$ws.Range("N5").Value = 10.627944
$ws.Range("O5").Value = 0.07742843230727542
$ws.Range("P5").Value = 0.05298650897942091
$ws.Range("Q5").Value = 1.516221460168
$ws.Range("R5").Value = 9.097328761007999
$ws.Range("S5").Value = 0.004881603075978495
$ws.Range("T5").Value = 0.003340621752392098
$ws.Range("I6").Value = 0.0630466474718978
$ws.Range("J6").Value = 0.06304664747189781
$ws.Range("M6").Value = 15.825637
$ws.Range("N6").Value = 47.476911
$ws.Range("O6").Value = 0.2305910274224278
$ws.Range("P6").Value = 0.2367001341949739
$ws.Range("Q6").Value = 4.515486803511333
$ws.Range("R6").Value = 40.639381231602
$ws.Range("S6").Value = 0.01453799121608452
$ws.Range("T6").Value = 0.01492314991714143
$ws.Range("G7").Value = 3.253975333333333
$ws.Range("H7").Value = 9.761925999999999
$ws.Range("I7").Value = 0.7190065996349845
$ws.Range("J7").Value = 0.7190065996349846
$ws.Range("M7").Value = 9.071155000000001
$ws.Range("N7").Value = 27.213465
$ws.Range("O7").Value = 0.1321733179750106
$ws.Range("P7").Value = 0.135675019324872
$ws.Range("Q7").Value = 29.51731461484333
$ws.Range("R7").Value = 265.65583153359
$ws.Range("S7").Value = 0.09503348791968595
$ws.Range("T7").Value = 0.09755123430018701
$ws.Range("G8").Value = 3.253975333333333
$ws.Range("H8").Value = 9.761925999999999
$ws.Range("I8").Value = 0.7190065996349845
$ws.Range("J8").Value = 0.7190065996349846
$ws.Range("N8").Value = 61.82726199999999
$ws.Range("O8").Value = 0.3002893736556623
$ws.Range("P8").Value = 0.3082450164524775
$ws.Range("Q8").Value = 67.06146182517909
$ws.Range("R8").Value = 603.5531564266118
$ws.Range("S8").Value = 0.2159100414586771
$ws.Range("T8").Value = 0.2216302011339257
$ws.Range("G9").Value = 3.253975333333333
$ws.Range("H9").Value = 9.761925999999999
$ws.Range("I9").Value = 0.7190065996349845
$ws.Range("J9").Value = 0.7190065996349846
$ws.Range("M9").Value = 17.81090666666666
$ws.Range("N9").Value = 53.43272
$ws.Range("O9").Value = 0.2595178486396241
$ws.Range("P9").Value = 0.2663933210482557
$ws.Range("Q9").Value = 57.95625095763554
$ws.Range("R9").Value = 521.6062586187199
$ws.Range("S9").Value = 0.1865950458949627
$ws.Range("T9").Value = 0.1915385559323771
$ws.Range("G10").Value = 3.253975333333333
$ws.Range("H10").Value = 9.761925999999999
$ws.Range("I10").Value = 0.7190065996349845
$ws.Range("J10").Value = 0.7190065996349846
$ws.Range("M10").Value = 5.313972
$ws.Range("N10").Value = 10.627944
$ws.Range("O10").Value = 0.07742843230727542
$ws.Range("P10").Value = 0.05298650897942091
$ws.Range("Q10").Value = 17.291533810024
$ws.Range("R10").Value = 103.749202860144
$ws.Range("S10").Value = 0.05567155382832167
$ws.Range("T10").Value = 0.03809764964782201
$ws.Range("G11").Value = 3.253975333333333
$ws.Range("H11").Value = 9.761925999999999
$ws.Range("I11").Value = 0.7190065996349845
$ws.Range("J11").Value = 0.7190065996349846
$ws.Range("M11").Value = 15.825637
$ws.Range("N11").Value = 47.476911
$ws.Range("O11").Value = 0.2305910274224278
$ws.Range("P11").Value = 0.2367001341949739
$ws.Range("Q11").Value = 51.49623243228733
$ws.Range("R11").Value = 463.466091890586
$ws.Range("S11").Value = 0.1657964705333372
$ws.Range("T11").Value = 0.1701889586206727
$ws.Range("G12").Value = 0.9863516666666667
$ws.Range("H12").Value = 2.959055
$ws.Range("I12").Value = 0.2179467528931175
$ws.Range("J12").Value = 0.2179467528931176
$ws.Range("M12").Value = 9.071155000000001
$ws.Range("N12").Value = 27.213465
$ws.Range("O12").Value = 0.1321733179750106
$ws.Range("P12").Value = 0.135675019324872
$ws.Range("Q12").Value = 8.947348852841669
$ws.Range("R12").Value = 80.52613967557501
$ws.Range("S12").Value = 0.02880674547176309
$ws.Range("T12").Value = 0.02956992991056682
$ws.Range("G13").Value = 0.9863516666666667
$ws.Range("H13").Value = 2.959055
$ws.Range("I13").Value = 0.2179467528931175
$ws.Range("J13").Value = 0.2179467528931176
$ws.Range("N13").Value = 61.82726199999999
$ws.Range("O13").Value = 0.3002893736556623
$ws.Range("P13").Value = 0.3082450164524775
$ws.Range("Q13").Value = 20.32780763971222
$ws.Range("R13").Value = 182.95026875741
$ws.Range("S13").Value = 0.06544709391655967
$ws.Range("T13").Value = 0.06718100043130307
$ws.Range("G14").Value = 0.9863516666666667
$ws.Range("H14").Value = 2.959055
$ws.Range("I14").Value = 0.2179467528931175
$ws.Range("J14").Value = 0.2179467528931176
$ws.Range("M14").Value = 17.81090666666666
$ws.Range("N14").Value = 53.43272
$ws.Range("O14").Value = 0.2595178486396241
$ws.Range("P14").Value = 0.2663933210482557
$ws.Range("Q14").Value = 17.56781747551111
$ws.Range("R14").Value = 158.1103572796
$ws.Range("S14").Value = 0.05656107242881363
$ws.Range("T14").Value = 0.05805955931488113
$ws.Range("G15").Value = 0.9863516666666667
$ws.Range("H15").Value = 2.959055
$ws.Range("I15").Value = 0.2179467528931175
$ws.Range("J15").Value = 0.2179467528931176
$ws.Range("M15").Value = 5.313972
$ws.Range("N15").Value = 10.627944
$ws.Range("O15").Value = 0.07742843230727542
$ws.Range("P15").Value = 0.05298650897942091
$ws.Range("Q15").Value = 5.241445138820001
$ws.Range("R15").Value = 31.44867083292
$ws.Range("S15").Value = 0.01687527540297524
$ws.Range("T15").Value = 0.01154823757920681
$ws.Range("G16").Value = 0.9863516666666667
$ws.Range("H16").Value = 2.959055
$ws.Range("I16").Value = 0.2179467528931175
$ws.Range("J16").Value = 0.2179467528931176
$ws.Range("M16").Value = 15.825637
$ws.Range("N16").Value = 47.476911
$ws.Range("O16").Value = 0.2305910274224278
$ws.Range("P16").Value = 0.2367001341949739
$ws.Range("Q16").Value = 15.60964343101167
$ws.Range("R16").Value = 140.486790879105
$ws.Range("S16").Value = 0.05025656567300595
$ws.Range("T16").Value = 0.05158802565715975
